$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# e089de36-871b-429d-af92-e5c21d483c0f.md row (row 3).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-24 06:46:40"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the e089de36-871b-429d-af92-e5c21d483c0f... row (row 3).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-24 06:46:35"
$wsZhCn.Range("K3").Value = "2016-08-24 06:47:04"

# de-de sheet: "Latest HO Xliff Generate Date" counterpart (Correspond
# Handoff Datetime, row 3) plus "Correspond Handback DateTime" (row 3).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-24 06:46:40"
$wsDeDe.Range("K3").Value = "2016-08-24 06:47:14"
